$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.213215351104736
$ws.Range("B1").Value = 1.4027179479599
$ws.Range("C1").Value = 1.789689421653748
$ws.Range("D1").Value = 3.83499813079834
$ws.Range("E1").Value = 3.460132360458374
